$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.097.69"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "1.661.64"
$ws.Range("E3").Value = "  -1.19%  "

$c = $ws.Range("D4")
$c.Value = "'1.002"
$c.Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  -0.22%  "

$c = $ws.Range("D5")
$c.Value = "'207.89"
$c.Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("E6").Value = "  -2.65%  "

$ws.Range("E7").Value = "  -0.21%  "

$c = $ws.Range("D8")
$c.Value = "'0.2580"
$c.Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  -4.08%  "

$c = $ws.Range("D9")
$c.Value = "'0.06297"
$c.Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  -0.30%  "

$c = $ws.Range("D10")
$c.Value = "'21.00"
$c.Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -1.69%  "

$c = $ws.Range("D11")
$c.Value = "'0.07533"
$c.Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "1.661.96"
$ws.Range("E12").Value = "  -1.19%  "

$ws.Range("E13").Value = "  -1.89%  "

$c = $ws.Range("D14")
$c.Value = "'0.5379"
$c.Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  -5.30%  "

$c = $ws.Range("D15")
$c.Value = "'66.18"
$c.Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").Value = "0.0₅7930"
$ws.Range("E16").Value = "  -2.86%  "

$ws.Range("D17").Value = "26.115.75"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("E18").Value = "  -0.16%  "

$c = $ws.Range("D19")
$c.Value = "'4.692"
$c.Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  -3.52%  "

$c = $ws.Range("D20")
$c.Value = "'187.50"
$c.Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("E21").Value = "  -4.04%  "

$c = $ws.Range("D22")
$c.Value = "'6.183"
$c.Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -0.76%  "

$c = $ws.Range("D23")
$c.Value = "'1.003"
$c.Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  -0.21%  "

$c = $ws.Range("D24")
$c.Value = "'148.30"
$c.Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  -4.15%  "

$c = $ws.Range("D26")
$c.Value = "'7.386"
$c.Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  -3.23%  "

$ws.Range("E27").Value = "  -2.37%  "

$ws.Range("E28").Value = "  +2.90%  "

$c = $ws.Range("D29")
$c.Value = "'0.06147"
$c.Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -4.94%  "

$c = $ws.Range("D30")
$c.Value = "'1.259"
$c.Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  -2.33%  "

$c = $ws.Range("D31")
$c.Value = "'3.467"
$c.Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  -2.12%  "

$c = $ws.Range("D32")
$c.Value = "'3.395"
$c.Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("E33").Value = "  -2.02%  "

$c = $ws.Range("D34")
$c.Value = "'0.9861"
$c.Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  -2.64%  "

$c = $ws.Range("D35")
$c.Value = "'2.386"
$c.Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  -1.27%  "

$c = $ws.Range("D36")
$c.Value = "'2.747"
$c.Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  +0.99%  "

$c = $ws.Range("D37")
$c.Value = "'0.5870"
$c.Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").Value = "1.104.01"
$ws.Range("E38").Value = "  -0.04%  "

$c = $ws.Range("D39")
$c.Value = "'0.01591"
$c.Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  -1.96%  "

$c = $ws.Range("D40")
$c.Value = "'5.978"
$c.Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  -3.51%  "

$c = $ws.Range("D41")
$c.Value = "'0.8452"
$c.Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("E42").Value = "  -0.47%  "

$c = $ws.Range("D43")
$c.Value = "'99.84"
$c.Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").Value = "1.811.12"
$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("D45").Value = "0.0₈108"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("E46").Value = "  -0.36%  "

$c = $ws.Range("D47")
$c.Value = "'54.90"
$c.Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  -4.09%  "

$c = $ws.Range("D48")
$c.Value = "'7.995"
$c.Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.32%  "

$c = $ws.Range("D49")
$c.Value = "'0.05235"
$c.Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  -0.67%  "

$c = $ws.Range("D50")
$c.Value = "'0.4242"
$c.Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  -0.76%  "

$c = $ws.Range("D51")
$c.Value = "'5.854"
$c.Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  -2.06%  "
